$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1158.8
$ws.Range("I2").Value = 1158.8
$ws.Range("K2").Value = 1158.8
$ws.Range("M2").Value = -1045.8

$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()

$ws.Range("H80").Value = 694.4
$ws.Range("I80").Value = 699.2857
$ws.Range("J80").Value = 683
$ws.Range("K80").Value = 2097.8571
$ws.Range("L80").Value = 2049
$ws.Range("M80").Value = -1099.8571
$ws.Range("N80").Value = -4045

$ws.Range("H83").Value = 694.4
$ws.Range("I83").Value = 699.2857
$ws.Range("J83").Value = 683
$ws.Range("K83").Value = 6293.571300000001
$ws.Range("L83").Value = 6147
$ws.Range("M83").Value = -1301.571300000001
$ws.Range("N83").Value = -16131

$ws.Range("H100").Value = 7106
$ws.Range("I100").Value = 2875.8
$ws.Range("K100").Value = 2875.8
$ws.Range("M100").Value = -2334.8

$ws.Range("H112").Value = 1788.9375
$ws.Range("I112").Value = 850
$ws.Range("J112").Value = 1923.0714
$ws.Range("K112").Value = 2550
$ws.Range("L112").Value = 5769.2142
$ws.Range("M112").Value = -1442
$ws.Range("N112").Value = -7985.2142

$ws.Range("H138").Value = 3916.35
$ws.Range("I138").Value = 2570.4285
$ws.Range("J138").Value = 4201.8486
$ws.Range("K138").Value = 7711.2855
$ws.Range("L138").Value = 12605.5458
$ws.Range("M138").Value = -2571.2855
$ws.Range("N138").Value = -22885.5458

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 20000
$ws.Range("J34").Value = 20000
$ws.Range("L34").Value = 20000
$ws.Range("N34").Value = -20542

$ws.Range("H37").Value = 37997
$ws.Range("J37").Value = 37997
$ws.Range("L37").Value = 37997
$ws.Range("N37").Value = -38543

$ws.Range("H61").Value = 1658.3334
$ws.Range("I61").Value = 979.46155
$ws.Range("K61").Value = 979.46155
$ws.Range("M61").Value = -767.46155

$ws.Range("H136").Value = 1658.3334
$ws.Range("I136").Value = 979.46155
$ws.Range("K136").Value = 2938.38465
$ws.Range("M136").Value = -388.38465

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2433.1667
$ws.Range("I86").Value = 2433.1667
$ws.Range("K86").Value = 2433.1667
$ws.Range("M86").Value = -1310.1667

$ws.Range("H89").Value = 2433.1667
$ws.Range("I89").Value = 2433.1667
$ws.Range("K89").Value = 12165.8335
$ws.Range("M89").Value = -6549.833500000001

$ws.Range("H99").Value = 996.5
$ws.Range("J99").Value = 993
$ws.Range("L99").Value = 993
$ws.Range("N99").Value = -3989

$ws.Range("H105").Value = 6124.25
$ws.Range("J105").Value = 7500
$ws.Range("L105").Value = 7500
$ws.Range("N105").Value = -10994

$ws.Range("H134").Value = 2214.5715
$ws.Range("I134").Value = 2059.1765
$ws.Range("K134").Value = 6177.529500000001
$ws.Range("M134").Value = -3642.529500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1890.6
$ws.Range("I16").Value = 1726.5
$ws.Range("K16").Value = 1726.5
$ws.Range("M16").Value = -1439.5

$ws.Range("H31").Value = 3828.8667
$ws.Range("I31").Value = 2563.6667
$ws.Range("K31").Value = 2563.6667
$ws.Range("M31").Value = -2268.6667

$ws.Range("H34").Value = 3828.8667
$ws.Range("I34").Value = 2563.6667
$ws.Range("K34").Value = 2563.6667
$ws.Range("M34").Value = -2361.6667

$ws.Range("H47").Value = 1001
$ws.Range("I47").Value = 1001
$ws.Range("K47").Value = 1001
$ws.Range("M47").Value = -435

$ws.Range("H58").Value = 2564.2856
$ws.Range("I58").Value = 2185.5454
$ws.Range("K58").Value = 2185.5454
$ws.Range("M58").Value = -1982.5454

$ws.Range("H60").Value = 35777.77
$ws.Range("J60").Value = 36808.88
$ws.Range("L60").Value = 36808.88
$ws.Range("N60").Value = -37830.88

$ws.Range("H62").Value = 1999.6666
$ws.Range("I62").Value = 1999.6666
$ws.Range("K62").Value = 1999.6666
$ws.Range("M62").Value = -1375.6666

$ws.Range("H65").Value = 1999.6666
$ws.Range("I65").Value = 1999.6666
$ws.Range("K65").Value = 9998.333000000001
$ws.Range("M65").Value = -6878.333000000001

$ws.Range("H86").Value = 19236.182
$ws.Range("I86").Value = 9133
$ws.Range("J86").Value = 29339.363
$ws.Range("K86").Value = 9133
$ws.Range("L86").Value = 29339.363
$ws.Range("M86").Value = -8010
$ws.Range("N86").Value = -31585.363

$ws.Range("H89").Value = 19236.182
$ws.Range("I89").Value = 9133
$ws.Range("J89").Value = 29339.363
$ws.Range("K89").Value = 45665
$ws.Range("L89").Value = 146696.815
$ws.Range("M89").Value = -40049
$ws.Range("N89").Value = -157928.815

$ws.Range("H107").Value = 2214.0667
$ws.Range("I107").Value = 1273
$ws.Range("J107").Value = 2684.6
$ws.Range("K107").Value = 1273
$ws.Range("L107").Value = 2684.6
$ws.Range("M107").Value = 647
$ws.Range("N107").Value = -6524.6

$ws.Range("H113").Value = 1890.6
$ws.Range("I113").Value = 1726.5
$ws.Range("K113").Value = 1726.5
$ws.Range("M113").Value = 443.5

$ws.Range("H122").Value = 2037.4286
$ws.Range("I122").Value = 670.6667
$ws.Range("K122").Value = 2012.0001
$ws.Range("M122").Value = 437.9999

$ws.Range("H132").Value = 3057.0833
$ws.Range("I132").Value = 2631.889
$ws.Range("J132").Value = 4332.6665
$ws.Range("K132").Value = 7895.667
$ws.Range("L132").Value = 12997.9995
$ws.Range("M132").Value = -5365.667
$ws.Range("N132").Value = -18057.9995

$ws.Range("H134").Value = 2548.1667
$ws.Range("I134").Value = 2231.5557
$ws.Range("K134").Value = 6694.6671
$ws.Range("M134").Value = -4159.6671

$ws.Range("H136").Value = 2564.2856
$ws.Range("I136").Value = 2185.5454
$ws.Range("K136").Value = 6556.6362
$ws.Range("M136").Value = -4006.6362

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 290
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

$ws.Range("H34").Value = 4654.222
$ws.Range("J34").Value = 6298.3335
$ws.Range("L34").Value = 18895.0005
$ws.Range("N34").Value = -19063.0005

$ws.Range("H55").Value = 5864.8335
$ws.Range("I55").Value = 1422.25
$ws.Range("K55").Value = 4266.75
$ws.Range("M55").Value = -4089.75

$ws.Range("H92").Value = 1500
$ws.Range("I92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("M92").ClearContents()

$ws.Range("H131").Value = 2192.75
$ws.Range("J131").Value = 2192.75
$ws.Range("L131").Value = 6578.25
$ws.Range("N131").Value = -16658.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 10001.333
$ws.Range("I80").Value = 9999
$ws.Range("K80").Value = 9999
$ws.Range("M80").Value = -9001

$ws.Range("H83").Value = 10001.333
$ws.Range("I83").Value = 9999
$ws.Range("K83").Value = 49995
$ws.Range("M83").Value = -45003

$ws.Range("H102").Value = 2500
$ws.Range("I102").Value = 2000
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 2000
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = -378
$ws.Range("N102").Value = -6244

$ws.Range("H122").Value = 1161.75
$ws.Range("I122").Value = 1161.75
$ws.Range("K122").Value = 3485.25
$ws.Range("M122").Value = -1035.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 417.42856
$ws.Range("J55").Value = 600
$ws.Range("L55").Value = 600
$ws.Range("N55").Value = -946

$ws.Range("H122").Value = 8328.286
$ws.Range("I122").Value = 6899.5
$ws.Range("J122").Value = 8899.799999999999
$ws.Range("K122").Value = 20698.5
$ws.Range("L122").Value = 26699.4
$ws.Range("M122").Value = -18248.5
$ws.Range("N122").Value = -31599.4

$ws.Range("H132").Value = 4688.3
$ws.Range("I132").Value = 4127
$ws.Range("J132").Value = 5998
$ws.Range("K132").Value = 12381
$ws.Range("L132").Value = 17994
$ws.Range("M132").Value = -9851
$ws.Range("N132").Value = -23054

$ws.Range("H136").Value = 3132.111
$ws.Range("I136").Value = 3226
$ws.Range("K136").Value = 9678
$ws.Range("M136").Value = -7128

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2738.3872
$ws.Range("I132").Value = 2038.4762
$ws.Range("K132").Value = 6115.4286
$ws.Range("M132").Value = -3585.4286

$ws.Range("H136").Value = 2473.7
$ws.Range("I136").Value = 2149.1738
$ws.Range("K136").Value = 6447.5214
$ws.Range("M136").Value = -3897.5214
